# Ironskillet 10.0 backdrop changes
# - bump the iron-skillet-version tag comment from 0.0.1 -> 0.0.2
# - add credential-enforcement / real-time-detection "set" lines for the
#   Outbound-URL, Alert-Only-URL and Exception-URL url-filtering profiles

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("set commands")

# 1) Bump the version comment in place (row is unaffected by the later
#    row-insertions because it sits well above them).
$ws.Range("A163").Value = 'set tag iron-skillet-version comments ""version 0.0.2 for 10.0: version of this IronSkillet template file""'

# 2) Insert the two new "real-time-detection" lines right after each
#    "...Javascript Exploit Detection... mlav-policy-action ..." line.
#    Insertions are performed top-to-bottom, so later anchor rows are
#    offset by the rows already inserted above them.

# -- Outbound-URL block (original anchor row 341) --
$ws.Rows.Item(342).Insert()
$ws.Rows.Item(343).Insert()
$ws.Range("A342").Value = "set profiles url-filtering Outbound-URL credential-enforcement alert real-time-detection"
$ws.Range("A343").Value = "set profiles url-filtering Outbound-URL alert real-time-detection"

# -- Alert-Only-URL block (original anchor row 347, now shifted +2 -> 349) --
$ws.Rows.Item(350).Insert()
$ws.Rows.Item(351).Insert()
$ws.Range("A350").Value = "set profiles url-filtering Alert-Only-URL credential-enforcement alert real-time-detection"
$ws.Range("A351").Value = "set profiles url-filtering Alert-Only-URL alert real-time-detection"

# -- Exception-URL block (original anchor row 357, now shifted +4 -> 361) --
$ws.Rows.Item(362).Insert()
$ws.Rows.Item(363).Insert()
$ws.Range("A362").Value = "set profiles url-filtering Exception-URL credential-enforcement alert real-time-detection"
$ws.Range("A363").Value = "set profiles url-filtering Exception-URL alert real-time-detection"
